$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193. This shifts the existing rows
# 193..234 down to 194..235, carrying along their values and styles
# (the D column keeps its date-format style automatically).
$ws.Rows.Item(193).Insert()

# Populate the newly inserted (now blank) row 193 with the new record.
$ws.Range("A193").Value = 7
$ws.Range("B193").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C193").Value = "Ñuble"
$ws.Range("D193").Value = 44785
$ws.Range("E193").Value = 16
$ws.Range("F193").Value = 100112032
$ws.Range("G193").Value = "Zapallo italiano"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 100
$ws.Range("K193").Value = 20000
$ws.Range("L193").Value = 21000
$ws.Range("M193").Value = 20500
$ws.Range("N193").Value = "$/caja 50 unidades"
$ws.Range("O193").Value = "Región de Arica y Parinacota"
$ws.Range("P193").Value = 410
$ws.Range("Q193").Value = 50
$ws.Range("R193").Value = "Hortaliza"
